$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The REPSWITCH1_Fam picture references moved folders — rename them to the
# new Experimental_scripts/Pictures_Fam location. Only column A (the
# "image" column) holds these picture filenames.
for ($r = 2; $r -le 46; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value2
    if ($v -like "REPSWITCH1_Fam/*") {
        $cell.Value2 = $v -replace "^REPSWITCH1_Fam/", "Experimental_scripts/Pictures_Fam/"
    }
}

# Move the active selection to A4.
$ws.Range("A4").Select()
